$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1:E4").Select()
$ws.Range("E1:E4").ClearContents()
